$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "year"/"rate" data (columns A:B, rows 1-41) down by one
# row, making room for a new header row at the top. Using Copy(Destination)
# so the destination cells inherit the copied content *and* its underlying
# type (shared-string text for column A, numbers for column B) exactly as
# stored in the source range, instead of round-tripping through a generic
# Value property that would coerce numeric-looking text back into numbers.
$ws.Range("A1:B41").Copy($ws.Range("A2:B42"))

# Replace the old top-left values with the new header labels. Clear() first
# so the new header cells fall back to the default (unstyled) cell format
# instead of inheriting the data rows' style.
$ws.Range("A1").Clear()
$ws.Range("B1").Clear()
$ws.Range("A1").Value = "år"
$ws.Range("B1").Value = "rente"

# Restore the view to the top of the sheet with D6 selected (matching the
# saved workbook state), instead of the old scrolled-down A15 top-left cell
# with E44 selected.
[void]$ws.Range("D6").Select()
